# Updates the cryptocurrency "Price" (column D) and hourly-change
# "Volume(1h)" (column E) figures on the active sheet to reflect the latest
# scrape (GitHub Actions cron refresh).
#
# Many "Price" values look like plain numbers (e.g. "0.9993", "24.51") even
# though the source data stores them as literal text - note some prices use
# a dotted-thousands format such as "29.392.66", which is not a valid number
# at all. Assigning a numeric-looking string straight to .Value lets Excel
# "smart" convert it into a real number, changing the stored cell type and
# losing the original text formatting. To keep these values as literal text,
# matching the original workbook, we momentarily force a Text number format
# before writing the value, then restore the cell to the default "Normal"
# style so no stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.392.66"
Set-TextValue "E2" "  -0.31%  "
Set-TextValue "D3" "1.848.58"
Set-TextValue "E3" "  -0.15%  "
Set-TextValue "D4" "0.9993"
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "E5" "  -0.99%  "
Set-TextValue "D6" "0.6325"
Set-TextValue "E6" "  -3.37%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "D8" "0.07586"
Set-TextValue "E8" "  +1.15%  "
Set-TextValue "D9" "0.2970"
Set-TextValue "E9" "  -0.53%  "
Set-TextValue "D10" "24.51"
Set-TextValue "E10" "  +0.30%  "
Set-TextValue "D11" "2.660.31"
Set-TextValue "E11" "  +43.45%  "
Set-TextValue "D12" "0.07725"
Set-TextValue "E12" "  +1.23%  "
Set-TextValue "E13" "  -0.65%  "
Set-TextValue "D14" "0.6851"
Set-TextValue "E14" "  +0.02%  "
Set-TextValue "D15" "82.87"
Set-TextValue "E15" "  -0.96%  "
Set-TextValue "D16" "0.000009955"
Set-TextValue "E16" "  +4.38%  "
Set-TextValue "D17" "6.184"
Set-TextValue "E17" "  +0.72%  "
Set-TextValue "D18" "29.421.35"
Set-TextValue "E18" "  -0.35%  "
Set-TextValue "D19" "231.21"
Set-TextValue "E19" "  -2.14%  "
Set-TextValue "D21" "0.9999"
Set-TextValue "D22" "7.600"
Set-TextValue "E22" "  -1.11%  "
Set-TextValue "E23" "  -0.08%  "
Set-TextValue "D24" "154.72"
Set-TextValue "E24" "  -1.42%  "
Set-TextValue "D25" "0.1397"
Set-TextValue "E25" "  -1.64%  "
Set-TextValue "E26" "  -0.54%  "
Set-TextValue "E27" "  -0.71%  "
Set-TextValue "D28" "1.471"
Set-TextValue "E28" "  -0.98%  "
Set-TextValue "D29" "0.05811"
Set-TextValue "E29" "  -3.50%  "
Set-TextValue "D30" "1.268"
Set-TextValue "E30" "  +1.56%  "
Set-TextValue "D31" "4.123"
Set-TextValue "E31" "  -0.33%  "
Set-TextValue "D32" "4.019"
Set-TextValue "E32" "  -1.31%  "
Set-TextValue "D33" "1.869"
Set-TextValue "E33" "  +0.92%  "
Set-TextValue "E34" "  -1.53%  "
Set-TextValue "D35" "0.7172"
Set-TextValue "E35" "  -0.78%  "
Set-TextValue "D36" "2.771.16"
Set-TextValue "E36" "  +37.36%  "
Set-TextValue "D37" "2.599"
Set-TextValue "E37" "  +0.15%  "
Set-TextValue "D38" "1.249.27"
Set-TextValue "E38" "  +4.12%  "
Set-TextValue "D39" "2.793"
Set-TextValue "E39" "  -0.30%  "
Set-TextValue "D40" "0.01808"
Set-TextValue "E40" "  +1.47%  "
Set-TextValue "D41" "0.9061"
Set-TextValue "E41" "  -0.15%  "
Set-TextValue "D42" "6.078"
Set-TextValue "E42" "  -2.61%  "
Set-TextValue "D43" "0.9998"
Set-TextValue "D44" "101.40"
Set-TextValue "E44" "  -0.44%  "
Set-TextValue "D45" "67.24"
Set-TextValue "E45" "  +1.41%  "
Set-TextValue "D46" "7.326"
Set-TextValue "E46" "  -0.95%  "
Set-TextValue "D47" "9.203"
Set-TextValue "E47" "  +1.58%  "
Set-TextValue "D48" "0.4013"
Set-TextValue "E48" "  -1.05%  "
Set-TextValue "D49" "1.690"
Set-TextValue "E49" "  +2.22%  "
Set-TextValue "E50" "  -0.23%  "
Set-TextValue "E51" "  +0.09%  "
